$d = $word.ActiveDocument

# --- Change 1 ---------------------------------------------------------
# "...il mio periodo formativo all'interno dell'ufficio tecnico dell'azienda.
#  Una cosa positiva..."
#                              -> "...ufficio tecnico. Una cosa positiva..."
# (drops the trailing " dell'azienda" before the full stop)
$ok1 = $d.Content.Find.Execute(
    "all’interno dell’ufficio tecnico dell’azienda. Una cosa positiva",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "all’interno dell’ufficio tecnico. Una cosa positiva",
    2
)
if (-not $ok1) {
    Write-Output "WARNING: change 1 (ufficio tecnico) not applied"
}

# --- Change 2 ---------------------------------------------------------
# "...documentandomi bene ho potuto realizzare tutte le cose chieste."
#        -> "...documentandomi bene ho potuto completare la maggior parte
#             dei task richiesti."
$ok2 = $d.Content.Find.Execute(
    "ho potuto realizzare tutte le cose chieste.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "ho potuto completare la maggior parte dei task richiesti.",
    2
)
if (-not $ok2) {
    Write-Output "WARNING: change 2 (completare la maggior parte) not applied"
}

Write-Output ("change1=" + $ok1 + " change2=" + $ok2)
